$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1242.1538
$ws.Range("I28").Value = 1062.1111
$ws.Range("K28").Value = 1062.1111
$ws.Range("M28").Value = -577.1111000000001
$ws.Range("H40").Value = 4643.405
$ws.Range("I40").Value = 6067
$ws.Range("J40").Value = 2745.2778
$ws.Range("K40").Value = 6067
$ws.Range("L40").Value = 2745.2778
$ws.Range("M40").Value = -5892
$ws.Range("N40").Value = -3095.2778
$ws.Range("H43").Value = 1747.125
$ws.Range("J43").Value = 1846.7142
$ws.Range("L43").Value = 1846.7142
$ws.Range("N43").Value = -1984.7142
$ws.Range("H58").Value = 3679.4
$ws.Range("J58").Value = 8998.5
$ws.Range("L58").Value = 26995.5
$ws.Range("N58").Value = -27295.5
$ws.Range("H62").Value = 78405.86
$ws.Range("I62").Value = 172280.5
$ws.Range("K62").Value = 172280.5
$ws.Range("M62").Value = -171656.5
$ws.Range("H65").Value = 78405.86
$ws.Range("I65").Value = 172280.5
$ws.Range("K65").Value = 861402.5
$ws.Range("M65").Value = -858282.5
$ws.Range("H74").Value = 5312.375
$ws.Range("I74").Value = 4750
$ws.Range("J74").Value = 6999.5
$ws.Range("K74").Value = 4750
$ws.Range("L74").Value = 6999.5
$ws.Range("M74").Value = -3814
$ws.Range("N74").Value = -8871.5
$ws.Range("H77").Value = 5312.375
$ws.Range("I77").Value = 4750
$ws.Range("J77").Value = 6999.5
$ws.Range("K77").Value = 23750
$ws.Range("L77").Value = 34997.5
$ws.Range("M77").Value = -19070
$ws.Range("N77").Value = -44357.5
$ws.Range("H97").Value = 1245.125
$ws.Range("J97").Value = 1243.6666
$ws.Range("L97").Value = 3730.9998
$ws.Range("N97").Value = -4722.9998
$ws.Range("H113").Value = 8640.846
$ws.Range("I113").Value = 27400
$ws.Range("J113").Value = 7077.5835
$ws.Range("K113").Value = 27400
$ws.Range("L113").Value = 7077.5835
$ws.Range("M113").Value = -24146
$ws.Range("N113").Value = -13585.5835
$ws.Range("H132").Value = 2938.0142
$ws.Range("I132").Value = 2913.25
$ws.Range("K132").Value = 8739.75
$ws.Range("M132").Value = -6209.75
$ws.Range("H135").Value = 630.5
$ws.Range("I135").Value = 461.1
$ws.Range("J135").Value = 1477.5
$ws.Range("K135").Value = 4149.900000000001
$ws.Range("L135").Value = 13297.5
$ws.Range("M135").Value = -1614.900000000001
$ws.Range("N135").Value = -18367.5
$ws.Range("H137").Value = 89998.25
$ws.Range("I137").Value = 127932.64
$ws.Range("K137").Value = 383797.92
$ws.Range("M137").Value = -381247.92
$ws.Range("H138").Value = 4298.354
$ws.Range("I138").Value = 4754.1113
$ws.Range("J138").Value = 4225.107
$ws.Range("K138").Value = 14262.3339
$ws.Range("L138").Value = 12675.321
$ws.Range("M138").Value = -9122.333899999998
$ws.Range("N138").Value = -22955.321
$ws.Range("H141").Value = 24447.46
$ws.Range("I141").Value = 24447.46
$ws.Range("K141").Value = 73342.38
$ws.Range("M141").Value = -68162.38

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1556.7646
$ws.Range("I2").Value = 1274
$ws.Range("K2").Value = 1274
$ws.Range("M2").Value = -1161
$ws.Range("H32").Value = 1657.9894
$ws.Range("I32").Value = 1424
$ws.Range("K32").Value = 1424
$ws.Range("M32").Value = -1137
$ws.Range("H48").Value = 242900
$ws.Range("J48").Value = 242900
$ws.Range("L48").Value = 242900
$ws.Range("N48").Value = -243668
$ws.Range("H97").Value = 13380.588
$ws.Range("I97").Value = 9730.5
$ws.Range("K97").Value = 9730.5
$ws.Range("M97").Value = -9234.5
$ws.Range("H102").Value = 2420.75
$ws.Range("I102").Value = 2420.75
$ws.Range("K102").Value = 2420.75
$ws.Range("M102").Value = -798.75
$ws.Range("H116").Value = 1556.7646
$ws.Range("I116").Value = 1274
$ws.Range("K116").Value = 1274
$ws.Range("M116").Value = 1020
$ws.Range("H122").Value = 4607.95
$ws.Range("I122").Value = 3985.8
$ws.Range("J122").Value = 4815.3335
$ws.Range("K122").Value = 11957.4
$ws.Range("L122").Value = 14446.0005
$ws.Range("M122").Value = -9507.400000000001
$ws.Range("N122").Value = -19346.0005
$ws.Range("H132").Value = 9927.842000000001
$ws.Range("I132").Value = 10242.333
$ws.Range("K132").Value = 30726.999
$ws.Range("M132").Value = -28196.999
$ws.Range("H139").Value = 190111
$ws.Range("J139").Value = 190111
$ws.Range("L139").Value = 190111
$ws.Range("N139").Value = -200391

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1556.7646
$ws.Range("I3").Value = 1274
$ws.Range("K3").Value = 1274
$ws.Range("M3").Value = -1160
$ws.Range("H20").Value = 3912.7
$ws.Range("J20").Value = 4164.8887
$ws.Range("L20").Value = 4164.8887
$ws.Range("N20").Value = -4658.8887
$ws.Range("H69").Value = 39500
$ws.Range("J69").Value = 39500
$ws.Range("L69").Value = 39500
$ws.Range("N69").Value = -41122
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H72").Value = 39500
$ws.Range("J72").Value = 39500
$ws.Range("L72").Value = 118500
$ws.Range("N72").Value = -126612
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H86").Value = 13114.9375
$ws.Range("I86").Value = 16809.727
$ws.Range("J86").Value = 4986.4
$ws.Range("K86").Value = 16809.727
$ws.Range("L86").Value = 4986.4
$ws.Range("M86").Value = -15686.727
$ws.Range("N86").Value = -7232.4
$ws.Range("H89").Value = 13114.9375
$ws.Range("I89").Value = 16809.727
$ws.Range("J89").Value = 4986.4
$ws.Range("K89").Value = 84048.63499999999
$ws.Range("L89").Value = 24932
$ws.Range("M89").Value = -78432.63499999999
$ws.Range("N89").Value = -36164
$ws.Range("H94").Value = 1549.7142
$ws.Range("I94").Value = 1634.6
$ws.Range("K94").Value = 1634.6
$ws.Range("M94").Value = -1183.6
$ws.Range("H99").Value = 3226.5667
$ws.Range("I99").Value = 2704.4783
$ws.Range("K99").Value = 2704.4783
$ws.Range("M99").Value = -1206.4783
$ws.Range("H105").Value = 113542.11
$ws.Range("J105").Value = 2992.5
$ws.Range("L105").Value = 2992.5
$ws.Range("N105").Value = -6486.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 451
$ws.Range("I10").Value = 329
$ws.Range("K10").Value = 329
$ws.Range("M10").Value = -190
$ws.Range("H16").Value = 1471.76
$ws.Range("I16").Value = 1471.76
$ws.Range("K16").Value = 1471.76
$ws.Range("M16").Value = -1184.76
$ws.Range("H62").Value = 2225
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2225
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H94").Value = 4384.6924
$ws.Range("J94").Value = 4969.857
$ws.Range("L94").Value = 4969.857
$ws.Range("N94").Value = -5871.857
$ws.Range("H104").Value = 26000
$ws.Range("J104").Value = 26000
$ws.Range("L104").Value = 26000
$ws.Range("N104").Value = -31242
$ws.Range("H113").Value = 1471.76
$ws.Range("I113").Value = 1471.76
$ws.Range("K113").Value = 1471.76
$ws.Range("M113").Value = 698.24
$ws.Range("H122").Value = 4096.615
$ws.Range("I122").Value = 2918.8572
$ws.Range("K122").Value = 8756.571599999999
$ws.Range("M122").Value = -6306.571599999999
$ws.Range("H132").Value = 32111.621
$ws.Range("I132").Value = 51720.477
$ws.Range("J132").Value = 6375
$ws.Range("K132").Value = 155161.431
$ws.Range("L132").Value = 19125
$ws.Range("M132").Value = -152631.431
$ws.Range("N132").Value = -24185

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 1996.1666
$ws.Range("I28").Value = 1238.5
$ws.Range("K28").Value = 3715.5
$ws.Range("M28").Value = -3483.5
$ws.Range("H56").Value = 6699.722
$ws.Range("I56").Value = 6699.722
$ws.Range("K56").Value = 6699.722
$ws.Range("M56").Value = -6169.722
$ws.Range("H63").Value = 6399.7144
$ws.Range("I63").Value = 6399.7144
$ws.Range("K63").Value = 19199.1432
$ws.Range("M63").Value = -18450.1432
$ws.Range("H66").Value = 6399.7144
$ws.Range("I66").Value = 6399.7144
$ws.Range("K66").Value = 57597.4296
$ws.Range("M66").Value = -53853.4296
$ws.Range("H68").Value = 3455.451
$ws.Range("J68").Value = 3653.641
$ws.Range("L68").Value = 10960.923
$ws.Range("N68").Value = -12582.923
$ws.Range("H71").Value = 3455.451
$ws.Range("J71").Value = 3653.641
$ws.Range("L71").Value = 32882.769
$ws.Range("N71").Value = -40994.769
$ws.Range("H87").Value = 21024.477
$ws.Range("I87").Value = 17902.8
$ws.Range("K87").Value = 53708.39999999999
$ws.Range("M87").Value = -52460.39999999999
$ws.Range("H90").Value = 21024.477
$ws.Range("I90").Value = 17902.8
$ws.Range("K90").Value = 161125.2
$ws.Range("M90").Value = -154885.2
$ws.Range("H92").Value = 383.83334
$ws.Range("J92").Value = 353.25
$ws.Range("L92").Value = 1059.75
$ws.Range("N92").Value = -3555.75
$ws.Range("H114").Value = 3140
$ws.Range("J114").Value = 3750
$ws.Range("L114").Value = 11250
$ws.Range("N114").Value = -17758
$ws.Range("H118").Value = 1764.4445
$ws.Range("I118").Value = 1779
$ws.Range("K118").Value = 5337
$ws.Range("M118").Value = -4094
$ws.Range("H119").Value = 999
$ws.Range("I119").Value = 999
$ws.Range("K119").Value = 2997
$ws.Range("M119").Value = 1841
$ws.Range("H120").Value = 18461.46
$ws.Range("I120").Value = 9999.5
$ws.Range("K120").Value = 29998.5
$ws.Range("M120").Value = -25160.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 65319.645
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 65319.645
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 65319.645
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -66437.64499999999
$ws.Range("H51").Value = 95591.62
$ws.Range("J51").Value = 95591.62
$ws.Range("L51").Value = 95591.62
$ws.Range("N51").Value = -96609.62
$ws.Range("H80").Value = 5657.5
$ws.Range("I80").Value = 5189
$ws.Range("K80").Value = 5189
$ws.Range("M80").Value = -4191
$ws.Range("H83").Value = 5657.5
$ws.Range("I83").Value = 5189
$ws.Range("K83").Value = 25945
$ws.Range("M83").Value = -20953
$ws.Range("H97").Value = 1476.0526
$ws.Range("I97").Value = 1296.75
$ws.Range("K97").Value = 1296.75
$ws.Range("M97").Value = -800.75
$ws.Range("H102").Value = 24983.092
$ws.Range("I102").Value = 1596.7646
$ws.Range("K102").Value = 1596.7646
$ws.Range("M102").Value = 25.23540000000003
$ws.Range("H107").Value = 3733.8
$ws.Range("I107").Value = 3417.25
$ws.Range("K107").Value = 3417.25
$ws.Range("M107").Value = -1497.25
$ws.Range("H113").Value = 4860.3076
$ws.Range("I113").Value = 4532.3335
$ws.Range("J113").Value = 5141.4287
$ws.Range("K113").Value = 4532.3335
$ws.Range("L113").Value = 5141.4287
$ws.Range("M113").Value = -2362.3335
$ws.Range("N113").Value = -9481.4287
$ws.Range("H122").Value = 996508.4399999999
$ws.Range("I122").Value = 1489796.9
$ws.Range("J122").Value = 9931.666999999999
$ws.Range("K122").Value = 4469390.699999999
$ws.Range("L122").Value = 29795.001
$ws.Range("M122").Value = -4466940.699999999
$ws.Range("N122").Value = -34695.001
$ws.Range("H132").Value = 4376.5405
$ws.Range("I132").Value = 3828.3635
$ws.Range("J132").Value = 5180.533
$ws.Range("K132").Value = 11485.0905
$ws.Range("L132").Value = 15541.599
$ws.Range("M132").Value = -8955.0905
$ws.Range("N132").Value = -20601.599

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9108.933999999999
$ws.Range("I7").Value = 7367.6
$ws.Range("K7").Value = 7367.6
$ws.Range("M7").Value = -7255.6
$ws.Range("H10").Value = 172137.42
$ws.Range("J10").Value = 172137.42
$ws.Range("L10").Value = 172137.42
$ws.Range("N10").Value = -172417.42
$ws.Range("H61").Value = 4500
$ws.Range("I61").Value = 4500
$ws.Range("K61").Value = 4500
$ws.Range("M61").Value = -4298
$ws.Range("H93").Value = 4448.5
$ws.Range("I93").Value = 3400
$ws.Range("K93").Value = 3400
$ws.Range("M93").Value = -2152
$ws.Range("H100").Value = 3830.805
$ws.Range("I100").Value = 3198.8708
$ws.Range("K100").Value = 3198.8708
$ws.Range("M100").Value = -2657.8708
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 4500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -2330
$ws.Range("H122").Value = 14998.667
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 14998.667
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 44996.001
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -49896.001
$ws.Range("H126").Value = 9108.933999999999
$ws.Range("I126").Value = 7367.6
$ws.Range("K126").Value = 22102.8
$ws.Range("M126").Value = -19632.8
$ws.Range("H132").Value = 2500
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 7500
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -12560
$ws.Range("H136").Value = 37146.3
$ws.Range("I136").Value = 145941.28
$ws.Range("J136").Value = 4034.7827
$ws.Range("K136").Value = 437823.84
$ws.Range("L136").Value = 12104.3481
$ws.Range("M136").Value = -435273.84
$ws.Range("N136").Value = -17204.3481

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1407.1111
$ws.Range("I100").Value = 721.6667
$ws.Range("J100").Value = 1749.8334
$ws.Range("K100").Value = 1443.3334
$ws.Range("L100").Value = 3499.6668
$ws.Range("M100").Value = -902.3334
$ws.Range("N100").Value = -4581.6668
$ws.Range("H113").Value = 1398.7084
$ws.Range("I113").Value = 1080.5
$ws.Range("J113").Value = 1716.9166
$ws.Range("K113").Value = 3241.5
$ws.Range("L113").Value = 5150.7498
$ws.Range("M113").Value = -1071.5
$ws.Range("N113").Value = -9490.7498
$ws.Range("H126").Value = 2931.389
$ws.Range("I126").Value = 2616.25
$ws.Range("K126").Value = 7848.75
$ws.Range("M126").Value = -5378.75
$ws.Range("H132").Value = 173494.55
$ws.Range("I132").Value = 5243.6313
$ws.Range("J132").Value = 410292.16
$ws.Range("K132").Value = 15730.8939
$ws.Range("L132").Value = 1230876.48
$ws.Range("M132").Value = -13200.8939
$ws.Range("N132").Value = -1235936.48
$ws.Range("H136").Value = 7792.485
$ws.Range("I136").Value = 9450.799999999999
$ws.Range("K136").Value = 28352.4
$ws.Range("M136").Value = -25802.4
